$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 10
$ws.Range("A3").Value = 99
$ws.Range("A2").Select()
